# Junction_Flooding_251.xlsx edit:
#  - Row 5 values are updated to a rounded ("custom accuracy") representation.
#  - Row 6 (the extra simulation timestep) is removed entirely, shrinking the
#    used range from A1:AH6 down to A1:AH5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 5 with the new, rounded values (2 decimal places).
$ws.Range("B5").Value = 10.45
$ws.Range("C5").Value = 7.73
$ws.Range("D5").Value = 0.9
$ws.Range("E5").Value = 22.97
$ws.Range("F5").Value = 18.36
$ws.Range("G5").Value = 8.17
$ws.Range("H5").Value = 30.39
$ws.Range("I5").Value = 12.8
$ws.Range("J5").Value = 5.56
$ws.Range("K5").Value = 8.12
$ws.Range("L5").Value = 9.21
$ws.Range("O5").Value = 8.27
$ws.Range("P5").Value = 11.64
$ws.Range("R5").Value = 0.7
$ws.Range("S5").Value = 0.54
$ws.Range("T5").Value = 118.54
$ws.Range("U5").Value = 23
$ws.Range("V5").Value = 7.63
$ws.Range("W5").Value = 15.25
$ws.Range("X5").Value = 8.04
$ws.Range("Y5").Value = 1.34
$ws.Range("Z5").Value = 14.79
$ws.Range("AA5").Value = 6.74
$ws.Range("AB5").Value = 6.06
$ws.Range("AD5").Value = 9.68
$ws.Range("AG5").Value = 4.21
$ws.Range("AH5").Value = 9.54

# Remove row 6 (formerly the last data row) so the sheet now ends at row 5.
$ws.Rows(6).Delete()
